# Scheduled data refresh: update currentAveragePrice / Leve profit columns
# across all job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 372.15384
$ws.Cells.Item(2, 9).Value = 100.666664
$ws.Cells.Item(2, 10).Value = 604.8570999999999
$ws.Cells.Item(2, 11).Value = 100.666664
$ws.Cells.Item(2, 12).Value = 604.8570999999999
$ws.Cells.Item(2, 13).Value = 12.333336
$ws.Cells.Item(2, 14).Value = -830.8570999999999

$ws.Cells.Item(51, 8).Value = 152857.14
$ws.Cells.Item(51, 10).Value = 174999.83
$ws.Cells.Item(51, 12).Value = 174999.83
$ws.Cells.Item(51, 14).Value = -175967.83

$ws.Cells.Item(55, 8).Value = 368.9
$ws.Cells.Item(55, 9).Value = 386.25
$ws.Cells.Item(55, 10).Value = 299.5
$ws.Cells.Item(55, 11).Value = 386.25
$ws.Cells.Item(55, 12).Value = 299.5
$ws.Cells.Item(55, 13).Value = -172.25
$ws.Cells.Item(55, 14).Value = -727.5

$ws.Cells.Item(94, 8).Value = 3332.6667
$ws.Cells.Item(94, 9).Value = 3499.2727
$ws.Cells.Item(94, 11).Value = 3499.2727
$ws.Cells.Item(94, 13).Value = -3048.2727

$ws.Cells.Item(104, 8).Value = 999
$ws.Cells.Item(104, 10).Value = 0
$ws.Cells.Item(104, 12).Value = 0
$ws.Cells.Item(104, 14).Value = ""

$ws.Cells.Item(125, 8).Value = 7984
$ws.Cells.Item(125, 9).Value = 7300
$ws.Cells.Item(125, 10).Value = 10036
$ws.Cells.Item(125, 11).Value = 65700
$ws.Cells.Item(125, 12).Value = 90324
$ws.Cells.Item(125, 13).Value = -63240
$ws.Cells.Item(125, 14).Value = -95244

$ws.Cells.Item(135, 8).Value = 979.6
$ws.Cells.Item(135, 9).Value = 979.6
$ws.Cells.Item(135, 11).Value = 8816.4
$ws.Cells.Item(135, 13).Value = -6281.4

$ws.Cells.Item(137, 8).Value = 2665.6667
$ws.Cells.Item(137, 9).Value = 2498.5
$ws.Cells.Item(137, 10).Value = 3000
$ws.Cells.Item(137, 11).Value = 7495.5
$ws.Cells.Item(137, 12).Value = 9000
$ws.Cells.Item(137, 13).Value = -4945.5
$ws.Cells.Item(137, 14).Value = -14100

$ws.Cells.Item(141, 8).Value = 5500
$ws.Cells.Item(141, 9).Value = 5000
$ws.Cells.Item(141, 10).Value = 6000
$ws.Cells.Item(141, 11).Value = 15000
$ws.Cells.Item(141, 12).Value = 18000
$ws.Cells.Item(141, 13).Value = -9820
$ws.Cells.Item(141, 14).Value = -28360

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 300
$ws.Cells.Item(4, 9).Value = 275
$ws.Cells.Item(4, 11).Value = 275
$ws.Cells.Item(4, 13).Value = -159

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(75, 8).Value = 49213.5
$ws.Cells.Item(75, 9).Value = 49213.5
$ws.Cells.Item(75, 11).Value = 49213.5
$ws.Cells.Item(75, 13).Value = -48277.5

$ws.Cells.Item(78, 8).Value = 49213.5
$ws.Cells.Item(78, 9).Value = 49213.5
$ws.Cells.Item(78, 11).Value = 147640.5
$ws.Cells.Item(78, 13).Value = -142960.5

$ws.Cells.Item(86, 8).Value = 947.8333
$ws.Cells.Item(86, 9).Value = 671.75
$ws.Cells.Item(86, 11).Value = 671.75
$ws.Cells.Item(86, 13).Value = 451.25

$ws.Cells.Item(89, 8).Value = 947.8333
$ws.Cells.Item(89, 9).Value = 671.75
$ws.Cells.Item(89, 11).Value = 3358.75
$ws.Cells.Item(89, 13).Value = 2257.25

$ws.Cells.Item(99, 8).Value = 931
$ws.Cells.Item(99, 9).Value = 778.2857
$ws.Cells.Item(99, 11).Value = 778.2857
$ws.Cells.Item(99, 13).Value = 719.7143

$ws.Cells.Item(105, 8).Value = 3197.8076
$ws.Cells.Item(105, 9).Value = 2635.9443
$ws.Cells.Item(105, 10).Value = 4462
$ws.Cells.Item(105, 11).Value = 2635.9443
$ws.Cells.Item(105, 12).Value = 4462
$ws.Cells.Item(105, 13).Value = -888.9443000000001
$ws.Cells.Item(105, 14).Value = -7956

$ws.Cells.Item(120, 8).Value = 29780.5
$ws.Cells.Item(120, 10).Value = 29780.5
$ws.Cells.Item(120, 12).Value = 29780.5
$ws.Cells.Item(120, 14).Value = -39456.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 12851.069
$ws.Cells.Item(99, 9).Value = 7114.2666
$ws.Cells.Item(99, 11).Value = 7114.2666
$ws.Cells.Item(99, 13).Value = -5616.2666

$ws.Cells.Item(106, 8).Value = 12333
$ws.Cells.Item(106, 10).Value = 12333
$ws.Cells.Item(106, 12).Value = 12333
$ws.Cells.Item(106, 14).Value = -14857

$ws.Cells.Item(126, 8).Value = 12851.069
$ws.Cells.Item(126, 9).Value = 7114.2666
$ws.Cells.Item(126, 11).Value = 21342.7998
$ws.Cells.Item(126, 13).Value = -18872.7998

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 13).Value = ""

$ws.Cells.Item(108, 8).Value = 9048.241
$ws.Cells.Item(108, 9).Value = 799.6667
$ws.Cells.Item(108, 11).Value = 2399.0001
$ws.Cells.Item(108, 13).Value = 480.9998999999998

$ws.Cells.Item(113, 8).Value = 683.3333
$ws.Cells.Item(113, 9).Value = 525
$ws.Cells.Item(113, 11).Value = 1575
$ws.Cells.Item(113, 13).Value = 595

$ws.Cells.Item(120, 8).Value = 8999.75
$ws.Cells.Item(120, 9).Value = 1999.5
$ws.Cells.Item(120, 10).Value = 16000
$ws.Cells.Item(120, 11).Value = 5998.5
$ws.Cells.Item(120, 12).Value = 48000
$ws.Cells.Item(120, 13).Value = -1160.5
$ws.Cells.Item(120, 14).Value = -57676

$ws.Cells.Item(130, 8).Value = 3397.5
$ws.Cells.Item(130, 10).Value = 4495
$ws.Cells.Item(130, 12).Value = 13485
$ws.Cells.Item(130, 14).Value = -23525

$ws.Cells.Item(138, 8).Value = 4646.625
$ws.Cells.Item(138, 9).Value = 2906.5715
$ws.Cells.Item(138, 11).Value = 8719.7145
$ws.Cells.Item(138, 13).Value = -3579.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 250.8
$ws.Cells.Item(2, 9).Value = 80.76922999999999
$ws.Cells.Item(2, 11).Value = 80.76922999999999
$ws.Cells.Item(2, 13).Value = 32.23077000000001

$ws.Cells.Item(11, 8).Value = 77571496
$ws.Cells.Item(11, 9).Value = 77571496
$ws.Cells.Item(11, 11).Value = 77571496
$ws.Cells.Item(11, 13).Value = -77571357

$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 13).Value = ""

$ws.Cells.Item(19, 8).Value = 750
$ws.Cells.Item(19, 9).Value = 500
$ws.Cells.Item(19, 10).Value = 1000
$ws.Cells.Item(19, 11).Value = 500
$ws.Cells.Item(19, 12).Value = 1000
$ws.Cells.Item(19, 13).Value = -212
$ws.Cells.Item(19, 14).Value = -1576

$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 14).Value = ""

$ws.Cells.Item(70, 8).Value = 17396.9
$ws.Cells.Item(70, 9).Value = 13162.667
$ws.Cells.Item(70, 11).Value = 13162.667
$ws.Cells.Item(70, 13).Value = -12892.667

$ws.Cells.Item(73, 8).Value = 17396.9
$ws.Cells.Item(73, 9).Value = 13162.667
$ws.Cells.Item(73, 11).Value = 13162.667
$ws.Cells.Item(73, 13).Value = -12226.667

$ws.Cells.Item(126, 8).Value = 5764.8
$ws.Cells.Item(126, 9).Value = 5633.1665
$ws.Cells.Item(126, 11).Value = 16899.4995
$ws.Cells.Item(126, 13).Value = -14429.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3363.0527
$ws.Cells.Item(46, 9).Value = 2393.2
$ws.Cells.Item(46, 10).Value = 7000
$ws.Cells.Item(46, 11).Value = 2393.2
$ws.Cells.Item(46, 12).Value = 7000
$ws.Cells.Item(46, 13).Value = -2205.2
$ws.Cells.Item(46, 14).Value = -7376

$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).Value = ""

$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).Value = ""

$ws.Cells.Item(82, 8).Value = 202577.8
$ws.Cells.Item(82, 9).Value = 1950
$ws.Cells.Item(82, 10).Value = 336329.66
$ws.Cells.Item(82, 11).Value = 1950
$ws.Cells.Item(82, 12).Value = 336329.66
$ws.Cells.Item(82, 13).Value = -1589
$ws.Cells.Item(82, 14).Value = -337051.66

$ws.Cells.Item(85, 8).Value = 202577.8
$ws.Cells.Item(85, 9).Value = 1950
$ws.Cells.Item(85, 10).Value = 336329.66
$ws.Cells.Item(85, 11).Value = 1950
$ws.Cells.Item(85, 12).Value = 336329.66
$ws.Cells.Item(85, 13).Value = -702
$ws.Cells.Item(85, 14).Value = -338825.66

$ws.Cells.Item(122, 8).Value = 5563.905
$ws.Cells.Item(122, 9).Value = 2771.8572
$ws.Cells.Item(122, 10).Value = 6959.9287
$ws.Cells.Item(122, 11).Value = 8315.571599999999
$ws.Cells.Item(122, 12).Value = 20879.7861
$ws.Cells.Item(122, 13).Value = -5865.571599999999
$ws.Cells.Item(122, 14).Value = -25779.7861

$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(131, 14).Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 14).Value = ""

$ws.Cells.Item(54, 8).Value = 29862.25
$ws.Cells.Item(54, 10).Value = 29632.666
$ws.Cells.Item(54, 12).Value = 29632.666
$ws.Cells.Item(54, 14).Value = -30672.666

$ws.Cells.Item(104, 8).Value = 19774.25
$ws.Cells.Item(104, 10).Value = 19774.25
$ws.Cells.Item(104, 12).Value = 19774.25
$ws.Cells.Item(104, 14).Value = -26762.25

$ws.Cells.Item(136, 8).Value = 4055.4167
$ws.Cells.Item(136, 9).Value = 3819.6
$ws.Cells.Item(136, 10).Value = 5234.5
$ws.Cells.Item(136, 11).Value = 11458.8
$ws.Cells.Item(136, 12).Value = 15703.5
$ws.Cells.Item(136, 13).Value = -8908.799999999999
$ws.Cells.Item(136, 14).Value = -20803.5
